$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers updated)
$wb.Worksheets.Item(1).Name = "GNG_TO-16511686508372188"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168652574372"
$wb.Worksheets.Item(3).Name = "RS_TO-16511686525763676"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651168652651368"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511686527304"

# Sheet 1 - GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511686508072188.csv"
$ws1.Range("B3").Value = "GNG_stims-1651168650820251.csv"
$ws1.Range("B4").Value = "go_stims-16511686508212187.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168650836221.csv"

# Sheet 2 - NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511686525633662.csv"
$ws2.Range("B3").Value = "TB-16511686521632276.csv"
$ws2.Range("B4").Value = "ZB-match_6-16511686512232296.csv"
$ws2.Range("B5").Value = "OB-1651168652038221.csv"
$ws2.Range("B6").Value = "ZB-match_1-16511686514532216.csv"
$ws2.Range("B7").Value = "OB-1651168652077227.csv"
$ws2.Range("B8").Value = "TB-1651168652461368.csv"
$ws2.Range("B9").Value = "ZB-match_7-16511686514122224.csv"
$ws2.Range("B10").Value = "OB-16511686518462212.csv"

# Sheet 4 - TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651168652604372.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168652578365.csv"
$ws4.Range("B4").Value = "MM_stims-1651168652635367.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686526053655.csv"
$ws4.Range("B6").Value = "MM_stims-1651168652650399.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686526363666.csv"

# Sheet 5 - vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511686527153726.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686526984007.csv"
$ws5.Range("B4").Value = "SAT_stims-16511686526563716.csv"
$ws5.Range("B5").Value = "SAT_stims-16511686526823661.csv"
